$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> FAPs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Ccl11"
$ws.Cells.Item(2, 3).Value = "Ackr4"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.7503183333333333
$ws.Cells.Item(2, 8).Value = 2.250955
$ws.Cells.Item(2, 9).Value = 0.001871730106429624
$ws.Cells.Item(2, 10).Value = 0.001871730106429624
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 1.890486333333333
$ws.Cells.Item(2, 14).Value = 5.671459
$ws.Cells.Item(2, 15).Value = 0.9442400689667343
$ws.Cells.Item(2, 16).Value = 0.9442400689667344
$ws.Cells.Item(2, 17).Value = 1.418466554816111
$ws.Cells.Item(2, 18).Value = 12.766198993345
$ws.Cells.Item(2, 19).Value = 0.001767362564782221
$ws.Cells.Item(2, 20).Value = 0.001767362564782221

# Row 3: ECs -> sCs
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Ccl11"
$ws.Cells.Item(3, 3).Value = "Ackr4"
$ws.Cells.Item(3, 4).Value = "sCs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.7503183333333333
$ws.Cells.Item(3, 8).Value = 2.250955
$ws.Cells.Item(3, 9).Value = 0.001871730106429624
$ws.Cells.Item(3, 10).Value = 0.001871730106429624
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.1116383333333333
$ws.Cells.Item(3, 14).Value = 0.334915
$ws.Cells.Item(3, 15).Value = 0.05575993103326565
$ws.Cells.Item(3, 16).Value = 0.05575993103326566
$ws.Cells.Item(3, 17).Value = 0.08376428820277777
$ws.Cells.Item(3, 18).Value = 0.753878593825
$ws.Cells.Item(3, 19).Value = 0.0001043675416474028
$ws.Cells.Item(3, 20).Value = 0.0001043675416474028

# Row 4: FAPs -> FAPs
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Ccl11"
$ws.Cells.Item(4, 3).Value = "Ackr4"
$ws.Cells.Item(4, 4).Value = "FAPs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 383.1307676666667
$ws.Cells.Item(4, 8).Value = 1149.392303
$ws.Cells.Item(4, 9).Value = 0.9557508602453543
$ws.Cells.Item(4, 10).Value = 0.9557508602453542
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 1.890486333333333
$ws.Cells.Item(4, 14).Value = 5.671459
$ws.Cells.Item(4, 15).Value = 0.9442400689667343
$ws.Cells.Item(4, 16).Value = 0.9442400689667344
$ws.Cells.Item(4, 17).Value = 724.3034801533421
$ws.Cells.Item(4, 18).Value = 6518.731321380078
$ws.Cells.Item(4, 19).Value = 0.9024582581930889
$ws.Cells.Item(4, 20).Value = 0.902458258193089

# Row 5: FAPs -> sCs
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Ccl11"
$ws.Cells.Item(5, 3).Value = "Ackr4"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 383.1307676666667
$ws.Cells.Item(5, 8).Value = 1149.392303
$ws.Cells.Item(5, 9).Value = 0.9557508602453543
$ws.Cells.Item(5, 10).Value = 0.9557508602453542
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.1116383333333333
$ws.Cells.Item(5, 14).Value = 0.334915
$ws.Cells.Item(5, 15).Value = 0.05575993103326565
$ws.Cells.Item(5, 16).Value = 0.05575993103326566
$ws.Cells.Item(5, 17).Value = 42.77208035102723
$ws.Cells.Item(5, 18).Value = 384.948723159245
$ws.Cells.Item(5, 19).Value = 0.05329260205226528
$ws.Cells.Item(5, 20).Value = 0.05329260205226528

# Row 6: sCs -> FAPs
$ws.Cells.Item(6, 1).Value = "sCs"
$ws.Cells.Item(6, 2).Value = "Ccl11"
$ws.Cells.Item(6, 3).Value = "Ackr4"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 16.98778433333333
$ws.Cells.Item(6, 8).Value = 50.963353
$ws.Cells.Item(6, 9).Value = 0.04237740964821621
$ws.Cells.Item(6, 10).Value = 0.0423774096482162
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 1.890486333333333
$ws.Cells.Item(6, 14).Value = 5.671459
$ws.Cells.Item(6, 15).Value = 0.9442400689667343
$ws.Cells.Item(6, 16).Value = 0.9442400689667344
$ws.Cells.Item(6, 17).Value = 32.11517411578078
$ws.Cells.Item(6, 18).Value = 289.036567042027
$ws.Cells.Item(6, 19).Value = 0.04001444820886323
$ws.Cells.Item(6, 20).Value = 0.04001444820886322

# Row 7: sCs -> sCs
$ws.Cells.Item(7, 1).Value = "sCs"
$ws.Cells.Item(7, 2).Value = "Ccl11"
$ws.Cells.Item(7, 3).Value = "Ackr4"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 16.98778433333333
$ws.Cells.Item(7, 8).Value = 50.963353
$ws.Cells.Item(7, 9).Value = 0.04237740964821621
$ws.Cells.Item(7, 10).Value = 0.0423774096482162
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.1116383333333333
$ws.Cells.Item(7, 14).Value = 0.334915
$ws.Cells.Item(7, 15).Value = 0.05575993103326565
$ws.Cells.Item(7, 16).Value = 0.05575993103326566
$ws.Cells.Item(7, 17).Value = 1.896487929999445
$ws.Cells.Item(7, 18).Value = 17.068391369995
$ws.Cells.Item(7, 19).Value = 0.002362961439352982
$ws.Cells.Item(7, 20).Value = 0.002362961439352982
